# Add "hydrogen combined cycle" as a new power plant type, and rename the
# existing "hydrogen" plant type to "hydrogen combustion turbine", across
# the three SoESCaOMCbIC sheets (capital, fixedOM, variableOM).

$wb = $excel.ActiveWorkbook

# Remember whichever sheet is active right now so we can restore it at the
# end (selecting ranges on other sheets below will otherwise leave one of
# them as the active/selected tab).
$originalActiveSheet = $wb.ActiveSheet

$sheetNames = @("SoESCaOMCbIC-capital", "SoESCaOMCbIC-fixedOM", "SoESCaOMCbIC-variableOM")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Rename row 24 ("hydrogen") to "hydrogen combustion turbine" - the
    # formulas already in B24:AQ24 (=B13 copied across) stay as-is.
    $ws.Range("A24").Value = "hydrogen combustion turbine"

    # Insert the new "hydrogen combined cycle" row directly below, mirroring
    # "natural gas combined cycle" (row 4) the same way row 24 mirrors
    # "natural gas peaker" (row 13).
    $ws.Range("A25").Value = "hydrogen combined cycle"
    $ws.Range("B25").Formula = "=B4"
    $ws.Range("C25:AQ25").Formula = "=C4"

    # Match the number formatting of the data row above.
    $ws.Range("B24:AQ24").Copy()
    $ws.Range("B25:AQ25").PasteSpecial(-4122)
    $ws.Application.CutCopyMode = $false

    # New label cells get a dedicated style: vertically centered, explicit
    # black font.
    $labelRange = $ws.Range("A24:A25")
    $labelRange.VerticalAlignment = -4108
    $labelRange.Font.Color = 0

    # Keep the on-screen selection in step with the newly-added row.
    $ws.Range("B25:AQ25").Select()
}

$originalActiveSheet.Activate()
